$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)  # Column E is the 5th column
    $cell.Value2 = $cell.Value2 - 1
}
